$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "I feel that the cooperation between the four of us worked",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "I feel that the cooperation among the four of us worked",
    2
)
